$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("D12").Value = "TensorFlow 2.6.0, Keras 2.6.0 Release"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/08/12/tensorflow-2-6-0-keras-2-6-0-release/"

# Row 16
$ws.Range("D16").Value = "Axiom-based Grad-CAM: Towards Accurate Visualization and Explanation of CNNs 내용 정리 [XAI-15]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/165"

# Row 21
$ws.Range("D21").Value = "[c++] 멀티 스레드(Multi Thread) _beginthreadex 사용법"
$ws.Range("E21").Value = "https://ms-review.tistory.com/24"

# Row 32
$ws.Range("D32").Value = "AdaBoost (에이다부스트)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/332"

# Row 36
$ws.Range("D36").Value = "Introduction to Autoencoder"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/330"

# Row 37
$ws.Range("D37").Value = "[Paper Review] Transferring inductive biases through knowledge distillation"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1805&mod=document&pageid=1"

# Row 52
$ws.Range("D52").Value = "숨은 DS"
